$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (ninth) test case was written for the "Read" row, bumping the
# Total Test Cases count for that row from 8 to 9. The Automated count
# (B8) is unchanged since the new test case isn't automated yet.
$ws.Range("C8").Value = 9

# G4 (=SUM($C:$C)) and G6 (=G5/G4) recalculate automatically from this edit.

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("B9").Select()
